# Added rows for new monitors (timer for delay requests) + a few label
# swaps on existing rows (the "browser" shared-string moved position).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# The big Selenium/WebDriver script body shared by every SCRIPT_API /
# SCRIPT_BROWSER monitor row.
$scriptBody = @'
var assert = require("assert");
$browser.get("http://example.com").then(function(){ 
  // Check the H1 title matches "Example Domain" 
  return $browser.findElement($driver.By.css("h1")).then(function(element){ 
    return element.getText().then(function(text){ 
      assert.equal("Example Domain", text, "Page H1 title did not match"); 
    }); 
  }); 
}).then(function(){ 
  // Check that the external link matches "http://www.iana.org/domains/example" 
  return $browser.findElement($driver.By.css("div > p > a")).then(function(element){ 
    return element.getAttribute("href").then(function(link){ 
      assert.equal("http://www.iana.org/domains/example", link, "More information link did not match"); 
    }); 
  }); 
});
'@

# ---------------------------------------------------------------------
# 1) Existing rows 2-5: two cells each swap which shared string they
#    point at ("browser"/"enabled" reshuffled, and a couple of the
#    monitor-name / type labels renumber) after the new strings were
#    inserted ahead of them in the shared-string table.
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 5).Value = "browser"
$ws.Cells.Item(2, 6).Value = "enabled"

$ws.Cells.Item(3, 1).Value = "monitor_testing_new"
$ws.Cells.Item(3, 5).Value = "SIMPLE"
$ws.Cells.Item(3, 6).Value = "enabled"

$ws.Cells.Item(4, 1).Value = "monitor_scripting_monitirng"
$ws.Cells.Item(4, 5).Value = "SCRIPT_BROWSER"
$ws.Cells.Item(4, 6).Value = "enabled"
$ws.Cells.Item(4, 12).Value = $scriptBody

$ws.Cells.Item(5, 1).Value = "monitor4"
$ws.Cells.Item(5, 5).Value = "SCRIPT_API"
$ws.Cells.Item(5, 6).Value = "enabled"
$ws.Cells.Item(5, 12).Value = $scriptBody

# Row 3 picks up an explicit custom height in the new layout.
$ws.Rows.Item(3).RowHeight = 39.75

# ---------------------------------------------------------------------
# 2) New rows 6-20: five repeating "monitor blocks" of varying shapes,
#    each cloned (formatting-wise) from whichever existing row has the
#    matching look, then the cell values are filled in.
# ---------------------------------------------------------------------

function New-BrowserRow($row, $name) {
    # Plain "browser" monitor row - same shape as row 2 (no script column).
    $ws.Range("A2:K2").Copy()
    $ws.Range("A$row`:K$row").PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = 15
    $ws.Cells.Item($row, 3).Value = "http://google.com"
    $ws.Cells.Item($row, 4).Value = "AWS_AP_SOUTH_1"
    $ws.Cells.Item($row, 5).Value = "browser"
    $ws.Cells.Item($row, 6).Value = "enabled"
    $ws.Cells.Item($row, 7).Value = 0.1
    $ws.Cells.Item($row, 9).Value = $false
    $ws.Cells.Item($row, 10).Value = $false
    $ws.Cells.Item($row, 11).Value = $false
}

function New-SimpleRow($row, $name) {
    # "SIMPLE" monitor row - same shape as row 3 (no script column).
    $ws.Range("A3:K3").Copy()
    $ws.Range("A$row`:K$row").PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = 15
    $ws.Cells.Item($row, 3).Value = "http://google.com"
    $ws.Cells.Item($row, 4).Value = "AWS_AP_SOUTH_1"
    $ws.Cells.Item($row, 5).Value = "SIMPLE"
    $ws.Cells.Item($row, 6).Value = "enabled"
    $ws.Cells.Item($row, 7).Value = 0.1
    $ws.Cells.Item($row, 9).Value = $false
    $ws.Cells.Item($row, 10).Value = $false
    $ws.Cells.Item($row, 11).Value = $false
    $ws.Rows.Item($row).RowHeight = 39.75
}

function New-ScriptApiRow($row, $name) {
    # "SCRIPT_API" monitor row with the big script body - same shape as row 5.
    $ws.Range("A5:L5").Copy()
    $ws.Range("A$row`:L$row").PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = 15
    $ws.Cells.Item($row, 3).Value = "http://google.com"
    $ws.Cells.Item($row, 4).Value = "AWS_AP_SOUTH_1"
    $ws.Cells.Item($row, 5).Value = "SCRIPT_API"
    $ws.Cells.Item($row, 6).Value = "enabled"
    $ws.Cells.Item($row, 7).Value = 0.1
    $ws.Cells.Item($row, 9).Value = $false
    $ws.Cells.Item($row, 10).Value = $false
    $ws.Cells.Item($row, 11).Value = $false
    $ws.Cells.Item($row, 12).Value = $scriptBody
    $ws.Rows.Item($row).RowHeight = 409.5
}

New-BrowserRow   6  "monitor112"
New-SimpleRow    7  "monitor_testing_new12"
New-ScriptApiRow 8  "monitorvv"
New-ScriptApiRow 9  "monitorvv4"
New-ScriptApiRow 10 "monitorvv3"
New-ScriptApiRow 11 "monitorvv2"
New-ScriptApiRow 12 "monitorvv1"
New-BrowserRow   13 "monitor1x"
New-SimpleRow    14 "monitor_testing_newx"
New-BrowserRow   15 "monitor1xx"
New-SimpleRow    16 "monitor_testing_newxx"
New-BrowserRow   17 "monitor1xx1"
New-SimpleRow    18 "monitor_testing_newxx1"
New-BrowserRow   19 "monitor1xx22"
New-SimpleRow    20 "monitor_testing_newxx12"

$excel.CutCopyMode = 0

# Rows 17-20 never had a script column to begin with (pasted from A2/A3
# which stop at L, but since they're plain/simple monitors, no script
# column data is used there) - nothing further to clear.

# ---------------------------------------------------------------------
# 3) Sheet cosmetics: wider name column, scrolled viewport, selection.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 35.67

$ws.Range("F21").Select()

Write-Output "done"
